# Weekly update: a new week of "Ajo" (garlic) price data was recorded for
# Terminal Hortofrutícola Agro Chillán. The new observation belongs at row 76
# (sheet is ordered so the new row is inserted there), so insert a blank row
# at 76, which shifts the former rows 76:181 down to 77:182 (and the sheet's
# used-range / dimension grows from R181 to R182 automatically), then fill
# the newly inserted row with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 76 - shifts existing rows 76:181 down to 77:182,
# carrying their values/styles with them (matches the diff: old row 76's
# data now lives at row 77, ..., old row 181's data now lives at row 182).
$ws.Rows("76:76").Insert()

# Populate the newly inserted row 76 with the new weekly observation. Most
# columns repeat the same catalog values as their neighboring rows for this
# market/product (Mercado, Región, Codreg, Categoría ID/Categoría, Variedad,
# Calidad, Unidad de comercialización, Origen, Kg o Unidades, Clasificación);
# only the date and the volume/price figures are new.
$ws.Range("A76").Value = 7
$ws.Range("B76").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C76").Value = "Ñuble"
$ws.Range("D76").Value = 44579
$ws.Range("E76").Value = 16
$ws.Range("F76").Value = 100112003
$ws.Range("G76").Value = "Ajo"
$ws.Range("H76").Value = "Chino"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 100
$ws.Range("K76").Value = 20000
$ws.Range("L76").Value = 21000
$ws.Range("M76").Value = 20500
$ws.Range("N76").Value = "$/caja 10 kilos"
$ws.Range("O76").Value = "China"
$ws.Range("P76").Value = 2050
$ws.Range("Q76").Value = 10
$ws.Range("R76").Value = "Hortaliza"
